$d = $word.ActiveDocument

# Paragraph 2 currently holds only the _GoBack bookmark (an empty
# paragraph). We rebuild it to:
#   "In f2" + " I " + " make " + "a new " + <bookmark _GoBack> + "change to f2.doc"
# as five distinct runs (matching the diff's run-per-insertion layout),
# and we append a brand-new empty paragraph right after it.
#
# InsertXML lets us inject literal OOXML runs/bookmark markup without
# Word's "merge adjacent identical-formatted runs" behaviour kicking in,
# and, when the injected fragment itself spans two <w:p> elements, the
# trailing one becomes a genuinely empty <w:p/> (no placeholder run) -
# exactly like the diff's new trailing paragraph.

$p2 = $d.Paragraphs.Item(2)
$insertPos = $p2.Range.Start
$target = $d.Range($insertPos, $insertPos)

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r><w:t>In f2</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> I </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> make </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">a new </w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t>change to f2.doc</w:t></w:r>' +
  '</w:p>' +
  '<w:p/>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($payload)

# The original bookmark paragraph got pushed down intact (it is now an
# empty paragraph holding only the stale "_GoBack" bookmark pair); the
# freshly rebuilt paragraph above already carries the live bookmark, so
# drop that now-redundant leftover paragraph.
$oldBookmarkPara = $d.Paragraphs.Item(4)
$oldBookmarkPara.Range.Delete()
